$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 7: Inscritos (E7) 13 -> 14
$ws.Range("E7").Value = 14

# Row 18: Inscritos (E18) 47 -> 48, Pagos (F18) 15 -> 16, Inscrições homologadas (H18) 15 -> 16
$ws.Range("E18").Value = 48
$ws.Range("F18").Value = 16
$ws.Range("H18").Value = 16
